$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (94 cell changes) ----
$ws = $wb.Worksheets.Item("ALC")
$edits = @(
    ,@("H12", 4132663.8)
    ,@("I12", 9091226)
    ,@("J12", 529)
    ,@("K12", 9091226)
    ,@("L12", 529)
    ,@("M12", -9091056)
    ,@("N12", -869)
    ,@("H15", 1048.075)
    ,@("I15", 1048.075)
    ,@("K15", 3144.225)
    ,@("M15", -2975.225)
    ,@("H19", 959.6)
    ,@("I19", 500)
    ,@("J19", 1074.5)
    ,@("K19", 500)
    ,@("L19", 1074.5)
    ,@("M19", -325)
    ,@("N19", -1424.5)
    ,@("H33", 7843097)
    ,@("I33", 2725896.8)
    ,@("K33", 2725896.8)
    ,@("M33", -2725667.8)
    ,@("H51", 76748.64)
    ,@("I51", 147113.58)
    ,@("J51", 6383.7144)
    ,@("K51", 147113.58)
    ,@("L51", 6383.7144)
    ,@("M51", -146629.58)
    ,@("N51", -7351.7144)
    ,@("H55", 1522.125)
    ,@("I55", 219.5)
    ,@("J55", 2824.75)
    ,@("K55", 219.5)
    ,@("L55", 2824.75)
    ,@("M55", -5.5)
    ,@("N55", -3252.75)
    ,@("H69", 14294980)
    ,@("J69", 10725)
    ,@("L69", 32175)
    ,@("N69", -33923)
    ,@("H72", 14294980)
    ,@("J72", 10725)
    ,@("L72", 96525)
    ,@("N72", -105261)
    ,@("I80", 306.8)
    ,@("K80", 920.4000000000001)
    ,@("M80", 77.59999999999991)
    ,@("I83", 306.8)
    ,@("K83", 2761.2)
    ,@("M83", 2230.8)
    ,@("H95", 34999)
    ,@("J95", 34999)
    ,@("L95", 34999)
    ,@("N95", -40491)
    ,@("H99", 83333864)
    ,@("I99", 588.8)
    ,@("J99", 500000260)
    ,@("K99", 1766.4)
    ,@("L99", 1500000780)
    ,@("M99", -268.3999999999999)
    ,@("N99", -1500003776)
    ,@("H116", 13774.889)
    ,@("I116", 13774.889)
    ,@("K116", 13774.889)
    ,@("M116", -10332.889)
    ,@("H125", 3789538.8)
    ,@("I125", 5683318)
    ,@("K125", 51149862)
    ,@("M125", -51147402)
    ,@("H132", 1762.0667)
    ,@("I132", 1703.5)
    ,@("K132", 5110.5)
    ,@("M132", -2580.5)
    ,@("H137", 1271897.4)
    ,@("I137", 5718.914)
    ,@("J137", 2279084.8)
    ,@("K137", 17156.742)
    ,@("L137", 6837254.399999999)
    ,@("M137", -14606.742)
    ,@("N137", -6842354.399999999)
    ,@("H138", 20131.572)
    ,@("I138", 59000)
    ,@("J138", 4584.2)
    ,@("K138", 177000)
    ,@("L138", 13752.6)
    ,@("M138", -171860)
    ,@("N138", -24032.6)
    ,@("H141", 4095)
    ,@("I141", 4095)
    ,@("J141", 0)
    ,@("K141", 12285)
    ,@("L141", 0)
    ,@("M141", -7105)
    ,@("N141", $null)
)
foreach ($e in $edits) {
    if ($null -eq $e[1]) {
        $ws.Range($e[0]).ClearContents()
    } else {
        $ws.Range($e[0]).Value = $e[1]
    }
}

# ---- Sheet: ARM (98 cell changes) ----
$ws = $wb.Worksheets.Item("ARM")
$edits = @(
    ,@("H2", 1316.0625)
    ,@("I2", 1284.5454)
    ,@("J2", 1385.4)
    ,@("K2", 1284.5454)
    ,@("L2", 1385.4)
    ,@("M2", -1171.5454)
    ,@("N2", -1611.4)
    ,@("H9", 35000)
    ,@("J9", 35000)
    ,@("L9", 35000)
    ,@("N9", -35340)
    ,@("H20", 35000)
    ,@("J20", 35000)
    ,@("L20", 35000)
    ,@("N20", -35540)
    ,@("H28", 16823.666)
    ,@("I28", 16823.666)
    ,@("K28", 16823.666)
    ,@("M28", -16631.666)
    ,@("H29", 25010350)
    ,@("I29", 50000450)
    ,@("J29", 20250)
    ,@("K29", 50000450)
    ,@("L29", 20250)
    ,@("M29", -50000142)
    ,@("N29", -20866)
    ,@("H30", 6225)
    ,@("I30", 3633.3333)
    ,@("J30", 14000)
    ,@("K30", 3633.3333)
    ,@("L30", 14000)
    ,@("M30", -3483.3333)
    ,@("N30", -14300)
    ,@("H32", 2248.6743)
    ,@("I32", 1188.3835)
    ,@("J32", 8202.615)
    ,@("K32", 1188.3835)
    ,@("L32", 8202.615)
    ,@("M32", -901.3834999999999)
    ,@("N32", -8776.615)
    ,@("H36", 14207.333)
    ,@("I36", 6311.25)
    ,@("J36", 29999.5)
    ,@("K36", 6311.25)
    ,@("L36", 29999.5)
    ,@("M36", -5965.25)
    ,@("N36", -30691.5)
    ,@("H61", 1765621.2)
    ,@("I61", 51049.91)
    ,@("J61", 4909002)
    ,@("K61", 51049.91)
    ,@("L61", 4909002)
    ,@("M61", -50837.91)
    ,@("N61", -4909426)
    ,@("H74", 573554.3)
    ,@("I74", 3409.6128)
    ,@("J74", 2782865)
    ,@("K74", 3409.6128)
    ,@("L74", 2782865)
    ,@("M74", -2535.6128)
    ,@("N74", -2784613)
    ,@("H77", 573554.3)
    ,@("I77", 3409.6128)
    ,@("J77", 2782865)
    ,@("K77", 17048.064)
    ,@("L77", 13914325)
    ,@("M77", -12680.064)
    ,@("N77", -13923061)
    ,@("H97", 4929.143)
    ,@("I97", 5498.1904)
    ,@("K97", 5498.1904)
    ,@("M97", -5002.1904)
    ,@("H99", 16823.666)
    ,@("I99", 16823.666)
    ,@("K99", 16823.666)
    ,@("M99", -13828.666)
    ,@("H102", 2695.0908)
    ,@("J102", 2475.4)
    ,@("L102", 2475.4)
    ,@("N102", -5719.4)
    ,@("H103", 69696)
    ,@("J103", 69696)
    ,@("L103", 69696)
    ,@("N103", -72040)
    ,@("H116", 1316.0625)
    ,@("I116", 1284.5454)
    ,@("J116", 1385.4)
    ,@("K116", 1284.5454)
    ,@("L116", 1385.4)
    ,@("M116", 1009.4546)
    ,@("N116", -5973.4)
    ,@("H136", 1765621.2)
    ,@("I136", 51049.91)
    ,@("J136", 4909002)
    ,@("K136", 153149.73)
    ,@("L136", 14727006)
    ,@("M136", -150599.73)
    ,@("N136", -14732106)
)
foreach ($e in $edits) {
    if ($null -eq $e[1]) {
        $ws.Range($e[0]).ClearContents()
    } else {
        $ws.Range($e[0]).Value = $e[1]
    }
}

# ---- Sheet: BSM (40 cell changes) ----
$ws = $wb.Worksheets.Item("BSM")
$edits = @(
    ,@("H3", 1316.0625)
    ,@("I3", 1284.5454)
    ,@("J3", 1385.4)
    ,@("K3", 1284.5454)
    ,@("L3", 1385.4)
    ,@("M3", -1170.5454)
    ,@("N3", -1613.4)
    ,@("H94", 1506)
    ,@("I94", 1590.3889)
    ,@("K94", 1590.3889)
    ,@("M94", -1139.3889)
    ,@("H104", 0)
    ,@("J104", 0)
    ,@("L104", 0)
    ,@("N104", $null)
    ,@("H105", 8266.41)
    ,@("I105", 5783.3076)
    ,@("J105", 13232.615)
    ,@("K105", 5783.3076)
    ,@("L105", 13232.615)
    ,@("M105", -4036.3076)
    ,@("N105", -16726.615)
    ,@("H107", 20343.824)
    ,@("I107", 22843)
    ,@("J107", 12221.5)
    ,@("K107", 22843)
    ,@("L107", 12221.5)
    ,@("M107", -20923)
    ,@("N107", -16061.5)
    ,@("H123", 0)
    ,@("J123", 0)
    ,@("L123", 0)
    ,@("N123", $null)
    ,@("H134", 29034784)
    ,@("I134", 2526.5518)
    ,@("J134", 450002500)
    ,@("K134", 7579.655400000001)
    ,@("L134", 1350007500)
    ,@("M134", -5044.655400000001)
    ,@("N134", -1350012570)
)
foreach ($e in $edits) {
    if ($null -eq $e[1]) {
        $ws.Range($e[0]).ClearContents()
    } else {
        $ws.Range($e[0]).Value = $e[1]
    }
}

# ---- Sheet: CRP (69 cell changes) ----
$ws = $wb.Worksheets.Item("CRP")
$edits = @(
    ,@("H7", 76.22221999999999)
    ,@("I7", 106)
    ,@("K7", 106)
    ,@("M7", 7)
    ,@("H22", 593.5833)
    ,@("I22", 607.5217)
    ,@("J22", 273)
    ,@("K22", 607.5217)
    ,@("L22", 273)
    ,@("M22", -257.5217)
    ,@("N22", -973)
    ,@("H31", 4849.9546)
    ,@("I31", 3235.8235)
    ,@("J31", 5236.4365)
    ,@("K31", 3235.8235)
    ,@("L31", 5236.4365)
    ,@("M31", -2940.8235)
    ,@("N31", -5826.4365)
    ,@("H32", 17000)
    ,@("I32", 4000)
    ,@("K32", 4000)
    ,@("M32", -3684)
    ,@("H33", 1119.75)
    ,@("I33", 1192.8)
    ,@("K33", 1192.8)
    ,@("M33", -813.8)
    ,@("H34", 4849.9546)
    ,@("I34", 3235.8235)
    ,@("J34", 5236.4365)
    ,@("K34", 3235.8235)
    ,@("L34", 5236.4365)
    ,@("M34", -3033.8235)
    ,@("N34", -5640.4365)
    ,@("H104", 66000)
    ,@("J104", 66000)
    ,@("L104", 66000)
    ,@("N104", -71242)
    ,@("H105", 2878.2)
    ,@("I105", 2263.5833)
    ,@("J105", 5336.6665)
    ,@("K105", 2263.5833)
    ,@("L105", 5336.6665)
    ,@("M105", -516.5832999999998)
    ,@("N105", -8830.666499999999)
    ,@("H106", 37266.2)
    ,@("J106", 37266.2)
    ,@("L106", 37266.2)
    ,@("N106", -39790.2)
    ,@("H107", 1595.1428)
    ,@("I107", 1653.8)
    ,@("J107", 1448.5)
    ,@("K107", 1653.8)
    ,@("L107", 1448.5)
    ,@("M107", 266.2)
    ,@("N107", -5288.5)
    ,@("H122", 3883.4546)
    ,@("I122", 2887.111)
    ,@("J122", 8367)
    ,@("K122", 8661.332999999999)
    ,@("L122", 25101)
    ,@("M122", -6211.332999999999)
    ,@("N122", -30001)
    ,@("H132", 2907.04)
    ,@("I132", 2979.2632)
    ,@("J132", 2678.3333)
    ,@("K132", 8937.7896)
    ,@("L132", 8034.999899999999)
    ,@("M132", -6407.7896)
    ,@("N132", -13094.9999)
)
foreach ($e in $edits) {
    if ($null -eq $e[1]) {
        $ws.Range($e[0]).ClearContents()
    } else {
        $ws.Range($e[0]).Value = $e[1]
    }
}

# ---- Sheet: CUL (36 cell changes) ----
$ws = $wb.Worksheets.Item("CUL")
$edits = @(
    ,@("H68", 2499.2856)
    ,@("J68", 3330)
    ,@("L68", 9990)
    ,@("N68", -11612)
    ,@("H71", 2499.2856)
    ,@("J71", 3330)
    ,@("L71", 29970)
    ,@("N71", -38082)
    ,@("H94", 434979.7)
    ,@("I94", 2109.5)
    ,@("J94", 513683.38)
    ,@("K94", 6328.5)
    ,@("L94", 1541050.14)
    ,@("M94", -5652.5)
    ,@("N94", -1542402.14)
    ,@("H104", 8470.888999999999)
    ,@("I104", 3047.6)
    ,@("J104", 15250)
    ,@("K104", 9142.799999999999)
    ,@("L104", 45750)
    ,@("M104", -6521.799999999999)
    ,@("N104", -50992)
    ,@("H109", 5550)
    ,@("I109", 660)
    ,@("J109", 30000)
    ,@("K109", 1980)
    ,@("L109", 90000)
    ,@("M109", -940)
    ,@("N109", -92080)
    ,@("H132", 1941.6471)
    ,@("I132", 1633.3334)
    ,@("J132", 2109.818)
    ,@("K132", 14700.0006)
    ,@("L132", 18988.362)
    ,@("M132", -12170.0006)
    ,@("N132", -24048.362)
)
foreach ($e in $edits) {
    if ($null -eq $e[1]) {
        $ws.Range($e[0]).ClearContents()
    } else {
        $ws.Range($e[0]).Value = $e[1]
    }
}

# ---- Sheet: GSM (62 cell changes) ----
$ws = $wb.Worksheets.Item("GSM")
$edits = @(
    ,@("H29", 668000)
    ,@("I29", 668000)
    ,@("J29", 0)
    ,@("K29", 668000)
    ,@("L29", 0)
    ,@("M29", -667710)
    ,@("N29", $null)
    ,@("H41", 5645.857)
    ,@("I41", 5645.857)
    ,@("K41", 5645.857)
    ,@("M41", -5290.857)
    ,@("H80", 19316438)
    ,@("I80", 41384.734)
    ,@("J80", 45600600)
    ,@("K80", 41384.734)
    ,@("L80", 45600600)
    ,@("M80", -40386.734)
    ,@("N80", -45602596)
    ,@("H83", 19316438)
    ,@("I83", 41384.734)
    ,@("J83", 45600600)
    ,@("K83", 206923.67)
    ,@("L83", 228003000)
    ,@("M83", -201931.67)
    ,@("N83", -228012984)
    ,@("H87", 0)
    ,@("J87", 0)
    ,@("L87", 0)
    ,@("N87", $null)
    ,@("H90", 0)
    ,@("J90", 0)
    ,@("L90", 0)
    ,@("N90", $null)
    ,@("H97", 405.875)
    ,@("I97", 437.125)
    ,@("J97", 343.375)
    ,@("K97", 437.125)
    ,@("L97", 343.375)
    ,@("M97", 58.875)
    ,@("N97", -1335.375)
    ,@("H99", 7197.5)
    ,@("I99", 4108.4443)
    ,@("K99", 4108.4443)
    ,@("M99", -1862.4443)
    ,@("H101", 0)
    ,@("J101", 0)
    ,@("L101", 0)
    ,@("N101", $null)
    ,@("H102", 33334842)
    ,@("I102", 45455460)
    ,@("J102", 3142.25)
    ,@("K102", 45455460)
    ,@("L102", 3142.25)
    ,@("M102", -45453838)
    ,@("N102", -6386.25)
    ,@("H107", 144549.28)
    ,@("I107", 200369)
    ,@("J107", 5000)
    ,@("K107", 200369)
    ,@("L107", 5000)
    ,@("M107", -198449)
    ,@("N107", -8840)
)
foreach ($e in $edits) {
    if ($null -eq $e[1]) {
        $ws.Range($e[0]).ClearContents()
    } else {
        $ws.Range($e[0]).Value = $e[1]
    }
}

# ---- Sheet: LTW (93 cell changes) ----
$ws = $wb.Worksheets.Item("LTW")
$edits = @(
    ,@("H4", 16981.5)
    ,@("I4", 16378.4)
    ,@("J4", 19997)
    ,@("K4", 16378.4)
    ,@("L4", 19997)
    ,@("M4", -16265.4)
    ,@("N4", -20223)
    ,@("H7", 11461.846)
    ,@("I7", 6501.375)
    ,@("K7", 6501.375)
    ,@("M7", -6389.375)
    ,@("H16", 8066.8335)
    ,@("J16", 7875.25)
    ,@("L16", 7875.25)
    ,@("N16", -8215.25)
    ,@("H28", 16981.5)
    ,@("I28", 16378.4)
    ,@("J28", 19997)
    ,@("K28", 16378.4)
    ,@("L28", 19997)
    ,@("M28", -16146.4)
    ,@("N28", -20461)
    ,@("H33", 12469.75)
    ,@("I33", 12469.75)
    ,@("K33", 12469.75)
    ,@("M33", -12179.75)
    ,@("H37", 16981.5)
    ,@("I37", 16378.4)
    ,@("J37", 19997)
    ,@("K37", 16378.4)
    ,@("L37", 19997)
    ,@("M37", -16271.4)
    ,@("N37", -20211)
    ,@("H40", 7405.6665)
    ,@("I40", 7870.857)
    ,@("J40", 5777.5)
    ,@("K40", 7870.857)
    ,@("L40", 5777.5)
    ,@("M40", -7734.857)
    ,@("N40", -6049.5)
    ,@("H46", 5514.36)
    ,@("I46", 22119.6)
    ,@("J46", 1363.05)
    ,@("K46", 22119.6)
    ,@("L46", 1363.05)
    ,@("M46", -21931.6)
    ,@("N46", -1739.05)
    ,@("H55", 622.53125)
    ,@("I55", 386.46667)
    ,@("J55", 830.82355)
    ,@("K55", 386.46667)
    ,@("L55", 830.82355)
    ,@("M55", -213.46667)
    ,@("N55", -1176.82355)
    ,@("H61", 8928.727999999999)
    ,@("I61", 4916.125)
    ,@("J61", 19629)
    ,@("K61", 4916.125)
    ,@("L61", 19629)
    ,@("M61", -4714.125)
    ,@("N61", -20033)
    ,@("H93", 3166.3333)
    ,@("I93", 3166.3333)
    ,@("J93", 0)
    ,@("K93", 3166.3333)
    ,@("L93", 0)
    ,@("M93", -1918.3333)
    ,@("N93", $null)
    ,@("H113", 8928.727999999999)
    ,@("I113", 4916.125)
    ,@("J113", 19629)
    ,@("K113", 4916.125)
    ,@("L113", 19629)
    ,@("M113", -2746.125)
    ,@("N113", -23969)
    ,@("H126", 11461.846)
    ,@("I126", 6501.375)
    ,@("K126", 19504.125)
    ,@("M126", -17034.125)
    ,@("H132", 4527.032)
    ,@("I132", 4132.2)
    ,@("J132", 5244.909)
    ,@("K132", 12396.6)
    ,@("L132", 15734.727)
    ,@("M132", -9866.599999999999)
    ,@("N132", -20794.727)
    ,@("H136", 49908.547)
    ,@("I136", 59599.11)
    ,@("J136", 6301)
    ,@("K136", 178797.33)
    ,@("L136", 18903)
    ,@("M136", -176247.33)
    ,@("N136", -24003)
)
foreach ($e in $edits) {
    if ($null -eq $e[1]) {
        $ws.Range($e[0]).ClearContents()
    } else {
        $ws.Range($e[0]).Value = $e[1]
    }
}

# ---- Sheet: WVR (54 cell changes) ----
$ws = $wb.Worksheets.Item("WVR")
$edits = @(
    ,@("H29", 0)
    ,@("I29", 0)
    ,@("K29", 0)
    ,@("M29", $null)
    ,@("H32", 10000)
    ,@("I32", 10000)
    ,@("K32", 10000)
    ,@("M32", -9683)
    ,@("H52", 6219.5)
    ,@("I52", 6219.5)
    ,@("K52", 6219.5)
    ,@("M52", -5993.5)
    ,@("H76", 111597.8)
    ,@("J76", 111597.8)
    ,@("L76", 111597.8)
    ,@("N76", -112227.8)
    ,@("H79", 111597.8)
    ,@("J79", 111597.8)
    ,@("L79", 111597.8)
    ,@("N79", -113781.8)
    ,@("H81", 172333.17)
    ,@("I81", 8000)
    ,@("K81", 16000)
    ,@("M81", -14939)
    ,@("H84", 172333.17)
    ,@("I84", 8000)
    ,@("K84", 80000)
    ,@("M84", -74696)
    ,@("H100", 1913.8572)
    ,@("I100", 1475.5)
    ,@("J100", 2498.3333)
    ,@("K100", 2951)
    ,@("L100", 4996.6666)
    ,@("M100", -2410)
    ,@("N100", -6078.6666)
    ,@("H112", 35999.668)
    ,@("J112", 35999.668)
    ,@("L112", 35999.668)
    ,@("N112", -38953.668)
    ,@("H122", 5014.615)
    ,@("I122", 4414)
    ,@("J122", 5715.3335)
    ,@("K122", 13242)
    ,@("L122", 17146.0005)
    ,@("M122", -10792)
    ,@("N122", -22046.0005)
    ,@("H126", 2184.6316)
    ,@("I126", 2184.6316)
    ,@("K126", 6553.8948)
    ,@("M126", -4083.8948)
    ,@("H136", 3834.5625)
    ,@("J136", 4639.857)
    ,@("L136", 13919.571)
    ,@("N136", -19019.571)
)
foreach ($e in $edits) {
    if ($null -eq $e[1]) {
        $ws.Range($e[0]).ClearContents()
    } else {
        $ws.Range($e[0]).Value = $e[1]
    }
}
